# Change in unit of AIC and run of new results
#
# Every worksheet in this workbook (tab per calendar year, 2000-2100) holds
# the same small A1:G8 table. The non-zero numeric results live in cells
# D5, E5, F5, G5 (row 5), D7, E7, F7, G7 (row 7) and D8, E8, F8, G8 (row 8).
# The unit of the underlying indicator (AIC) changed, which rescales every
# one of those numbers by a factor of 1e-6 (e.g. -2079.260599389473 becomes
# -0.002079260599389472). Zero-valued cells stay zero, so a blanket
# multiplication by 1e-6 is safe and reproduces the diff exactly.

$wb = $excel.ActiveWorkbook
$factor = 0.000001

foreach ($ws in $wb.Worksheets) {
    foreach ($rowNum in 5, 7, 8) {
        foreach ($col in 4, 5, 6, 7) {
            $cell = $ws.Cells.Item($rowNum, $col)
            $val = $cell.Value2
            if ($val -ne $null -and $val -ne 0) {
                $cell.Value2 = $val * $factor
            }
        }
    }
}
